# Add the "try_island" sheet and update P_req / Q_req data per commit.
$wb = $excel.ActiveWorkbook

# --- Update P_req (sheet1) ---
$wsP = $wb.Worksheets.Item("P_req")
for ($r = 2; $r -le 7; $r++) {
    $wsP.Cells.Item($r, 2).Value = -2000000
}
for ($r = 8; $r -le 52; $r++) {
    $wsP.Cells.Item($r, 2).Value = 0
}
$wsP.Range("B8:B11").Select() | Out-Null

# --- Update Q_req (sheet2) ---
$wsQ = $wb.Worksheets.Item("Q_req")
for ($r = 2; $r -le 52; $r++) {
    $wsQ.Cells.Item($r, 2).Value = 0
}
$wsQ.Range("B25").Select() | Out-Null

# --- Add try_island sheet (placed right after Q_req) ---
$wsNew = $wb.Worksheets.Add($null, $wsQ)
$wsNew.Name = "try_island"

$wsNew.Cells.Item(1, 1).Value = "Time [s]"
$wsNew.Cells.Item(1, 2).Value = "try_island"

for ($r = 2; $r -le 38; $r++) {
    $wsNew.Cells.Item($r, 1).Value = $r - 2
    $wsNew.Cells.Item($r, 2).Value = 0
}
for ($r = 39; $r -le 52; $r++) {
    $wsNew.Cells.Item($r, 1).Value = $r - 2
    $wsNew.Cells.Item($r, 2).Value = 1
}

$wsNew.Range("E16").Select() | Out-Null

# Restore P_req as the active/selected sheet, matching the original tab state.
$wsP.Activate() | Out-Null
